$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default/"Normal" data-row style) used to restore cell style
# after forcing text entry for numeric-looking Price values, so no stray
# quote-prefix style sticks to the cell (keeps D column plain "General"/default style).
$refStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "'28.329.43"
$ws.Range("D2").Style = $refStyle
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").Value = "'1.550.26"
$ws.Range("D3").Style = $refStyle

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'210.00"
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").Value = "'0.479"
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = "  -2.22%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "'23.83"
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = "  -0.71%  "

$ws.Range("E9").Value = "  -1.90%  "

$ws.Range("E10").Value = "  -1.53%  "

$ws.Range("D11").Value = "'0.0888"
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = "  -0.61%  "

$ws.Range("D12").Value = "'1.772.83"
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("D13").Value = "'1.553.97"
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = "  -1.52%  "

$ws.Range("D14").Value = "'28.330.58"
$ws.Range("D14").Style = $refStyle

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'3.62"
$ws.Range("D15").Style = $refStyle
$ws.Range("E15").Value = "  -1.90%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.509"
$ws.Range("D16").Style = $refStyle
$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("D17").Value = "'60.75"
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = "  -2.27%  "

$ws.Range("D18").Value = "'227.23"
$ws.Range("D18").Style = $refStyle
$ws.Range("E18").Value = "  -1.66%  "

$ws.Range("D19").Value = "'7.35"
$ws.Range("D19").Style = $refStyle
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").Value = "'3.91"
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("E23").Value = "  -2.49%  "

$ws.Range("E24").Value = "  -1.78%  "

$ws.Range("D25").Value = "'150.92"
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("E26").Value = "  -1.93%  "

$ws.Range("E27").Value = "  -1.18%  "

$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("D29").Value = "'6.24"
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = "  -3.09%  "

$ws.Range("E30").Value = "  -3.38%  "

$ws.Range("E31").Value = "  -4.86%  "

$ws.Range("D33").Value = "'1.382.50"
$ws.Range("D33").Style = $refStyle
$ws.Range("E33").Value = "  -1.15%  "

$ws.Range("D34").Value = "'3.00"
$ws.Range("D34").Style = $refStyle
$ws.Range("E34").Value = "  -3.05%  "

$ws.Range("E35").Value = "  +1.43%  "

$ws.Range("D36").Value = "'1.47"
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = "  -3.64%  "

$ws.Range("E37").Value = "  -1.01%  "

$ws.Range("D38").Value = "'2.57"
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = "  -1.80%  "

$ws.Range("E39").Value = "  -2.58%  "

$ws.Range("E40").Value = "  +1.71%  "

$ws.Range("D41").Value = "'0.509"
$ws.Range("D41").Style = $refStyle
$ws.Range("E41").Value = "  -2.32%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = $refStyle
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").Value = "'0.777"
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("E44").Value = "  -2.23%  "

$ws.Range("E45").Value = "  -1.71%  "

$ws.Range("D46").Value = "'61.93"
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = "  -1.79%  "

$ws.Range("D47").Value = "'1.685.07"
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = "  -1.71%  "

$ws.Range("D48").Value = "'0.874"
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = "  -9.15%  "

$ws.Range("D49").Value = "'85.23"
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = "  -1.55%  "

$ws.Range("D50").Value = "'42.13"
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = "  +5.41%  "

$ws.Range("E51").Value = "  +0.05%  "
